# "fixed corrupted xlsx file"
#
# - drops the stray "Matrix_Factorization" sheet
# - renames the "Best Runtime" header to "Runtime" and drops the
#   "Average Runtime" column entirely
# - the first data row's date is cleared out (left blank, but keeps its
#   date number-format), the label row now reads "11-12-2018" /
#   "11-13-2018" instead of real dates
# - the final runtime figures are tweaked slightly

$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook

function Fix-TimingSheet($ws, $finalRuntime) {
    # Header: "Best Runtime" -> "Runtime"
    $ws.Range("B1").Value = "Runtime"

    # Row 2: blank the date (format stays) and drop its runtime value
    $ws.Range("A2").ClearContents()
    $ws.Range("B2").Clear()

    # Row 3: swap the stored date serial for a plain text date label;
    # its runtime figure is unchanged
    $ws.Range("A3").Style = "Normal"
    $ws.Range("A3").Value = "'11-12-2018"
    $ws.Range("B3").Style = "Normal"

    # Row 4: updated runtime figure
    $ws.Range("B4").Value = $finalRuntime

    # The "Average Runtime" column is gone entirely
    $ws.Range("C1").EntireColumn.Delete()
}

$wsUser = $wb.Worksheets.Item("user_based")
Fix-TimingSheet $wsUser 1.53

$wsItem = $wb.Worksheets.Item("item_based")
Fix-TimingSheet $wsItem 4.56

# Remove the corrupt extra sheet entirely
$wb.Worksheets.Item("Matrix_Factorization").Delete() | Out-Null

Write-Output "timing workbook fixed"
